$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 0.1169995834814548
    "C2" = 1.626987699542094
    "D2" = 3.223369029078222
    "E2" = 13.86384647080068
    "G2" = 18.83120278290246

    "B3" = 0.04172184405617529
    "C3" = 0.04103571897497393
    "D3" = 0.1496068669990043
    "E3" = 0.5333859586016987
    "G3" = 0.7657503886318522

    "B4" = 1.445647641019636
    "C4" = 1.626987699542094
    "D4" = 0.1496068669990043
    "E4" = 0.5333859586016987
    "G4" = 3.755628166162433

    "B5" = 3.272327238179451
    "C5" = 1.626987699542094
    "D5" = 0.7210945179870265
    "E5" = 0.5333859586016987
    "G5" = 6.15379541431027
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
